$d = $word.ActiveDocument

# The document contains five "<id>...</id>" tags (p034r_1 .. p034r_5), each
# split across three runs: "<id>" (Courier New, gold), the bare id text
# (plain black), and "</id>" (Courier New, gold). Collapse each trio into a
# single run - formatted like the "<id>" run - whose text is the full
# "<id>p034r_N</id>" string, as described by the diff.
for ($i = 1; $i -le 5; $i++) {
    $tag = "<id>p034r_$i</id>"
    $rng = $d.Content
    $found = $rng.Find.Execute($tag, $true, $false, $false, $false, $false, $true, 0, $false, $tag, 2)
    Write-Host "id $i replaced:" $found
}
